# live_trading_results.xlsx update
# Commit: Trade #6 closed at 2026-02-17 20:02:27 - unknown UNKNOWN +0.000%

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Summary sheet
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B2").Value  = 1400            # Initial Capital
$wsSummary.Range("B3").Value  = 1399.76         # Current Capital
$wsSummary.Range("B4").Value  = -0.24           # Total P&L $
$wsSummary.Range("B5").Value  = -0.8            # Total P&L %
$wsSummary.Range("B6").Value  = 6               # Total Trades
$wsSummary.Range("B7").Value  = 2               # Winning Trades
$wsSummary.Range("B9").Value  = 33.33           # Win Rate %
$wsSummary.Range("B11").Value = 14              # Active Strategies

# ---------------------------------------------------------------------------
# 2) Strategy Status sheet
#    A new strategy "EMAArbitrage" is inserted (alphabetically) between
#    "CopyTrading" (row 2) and "HighProbConvergence" (row 3), shifting all
#    following rows down by one. The MarketMaking row (now row 5) is
#    refreshed with the latest stats.
# ---------------------------------------------------------------------------
$wsStrategy = $wb.Worksheets.Item("Strategy Status")

# Insert a new blank row above current row 3 (HighProbConvergence)
$wsStrategy.Rows.Item(3).Insert()

# New row 3: EMAArbitrage (fresh strategy, no trades yet)
$wsStrategy.Range("A3").Value = "EMAArbitrage"
$wsStrategy.Range("B3").Value = "ACTIVE"
$wsStrategy.Range("C3").Value = 100
$wsStrategy.Range("D3").Value = 0
$wsStrategy.Range("E3").Value = 0
$wsStrategy.Range("F3").Value = 0
$wsStrategy.Range("G3").Value = 0

# Row 4 is now HighProbConvergence - reset to default values
$wsStrategy.Range("C4").Value = 100
$wsStrategy.Range("D4").Value = 0
$wsStrategy.Range("E4").Value = 0
$wsStrategy.Range("F4").Value = 0
$wsStrategy.Range("G4").Value = 0

# Row 5 is now MarketMaking - refresh with latest trading stats
$wsStrategy.Range("C5").Value = 99.76000000000001
$wsStrategy.Range("D5").Value = 6
$wsStrategy.Range("E5").Value = -0.24
$wsStrategy.Range("F5").Value = -0.24
$wsStrategy.Range("G5").Value = 33.33

# ---------------------------------------------------------------------------
# 3) All Trades / 4) MarketMaking sheets (identical data)
#    Trades 1-6 get refreshed execution details, trades 7-11 are removed.
# ---------------------------------------------------------------------------
function Update-TradeSheet($ws) {
    # Row 2
    $ws.Range("C2").Value = "19:50:33"
    $ws.Range("E2").Value = "DOWN"
    $ws.Range("F2").Value = 0.49
    $ws.Range("G2").Value = 0.34
    $ws.Range("I2").Value = -30.6122
    $ws.Range("J2").Value = -0.15
    $ws.Range("K2").Value = 99.84999999999999
    $ws.Range("Q2").Value = 5.05

    # Row 3
    $ws.Range("C3").Value = "19:55:36"
    $ws.Range("F3").Value = 0.66
    $ws.Range("G3").Value = 0.59
    $ws.Range("I3").Value = -10.6061
    $ws.Range("J3").Value = -0.07000000000000001
    $ws.Range("K3").Value = 99.78

    # Row 4
    $ws.Range("C4").Value = "19:55:44"
    $ws.Range("F4").Value = 0.43
    $ws.Range("G4").Value = 0.42
    $ws.Range("I4").Value = -2.3256
    $ws.Range("J4").Value = -0.01
    $ws.Range("K4").Value = 99.77
    $ws.Range("Q4").Value = 0.13

    # Row 5
    $ws.Range("C5").Value = "19:55:52"
    $ws.Range("E5").Value = "UP"
    $ws.Range("F5").Value = 0.59
    $ws.Range("G5").Value = 0.6
    $ws.Range("I5").Value = 1.6949
    $ws.Range("J5").Value = 0.01
    $ws.Range("K5").Value = 99.78
    $ws.Range("Q5").Value = 0.14

    # Row 6
    $ws.Range("C6").Value = "19:56:00"
    $ws.Range("E6").Value = "UP"
    $ws.Range("F6").Value = 0.58
    $ws.Range("G6").Value = 0.55
    $ws.Range("I6").Value = -5.1724
    $ws.Range("J6").Value = -0.03
    $ws.Range("K6").Value = 99.75

    # Row 7
    $ws.Range("C7").Value = "19:56:08"
    $ws.Range("E7").Value = "DOWN"
    $ws.Range("F7").Value = 0.48
    $ws.Range("G7").Value = 0.49
    $ws.Range("I7").Value = 2.0833
    $ws.Range("J7").Value = 0.01
    $ws.Range("K7").Value = 99.76000000000001
    $ws.Range("Q7").Value = 6.34

    # Remove trades 7-11 (old rows 8-12)
    $ws.Range("A8:Q12").EntireRow.Delete()
}

$wsAllTrades = $wb.Worksheets.Item("All Trades")
Update-TradeSheet $wsAllTrades

$wsMarketMaking = $wb.Worksheets.Item("MarketMaking")
Update-TradeSheet $wsMarketMaking
